$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.080.41'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.397.57'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.75%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '568.76'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.18%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '156.11'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.77%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.609'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.399.53'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('E11').Value = '  -3.67%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.438'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.982.19'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('E15').Value = '  -4.75%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '27.73'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.144.28'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.405.46'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.93'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '372.79'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.79%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.93'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.92%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '71.75'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.26%  '
$ws.Range('E26').Value = '  -5.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.89'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.71%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.175'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.46'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.07'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.00'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '23.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.20'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.59'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.83%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '159.77'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.88'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0756'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.63%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '26.64'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.41%  '
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.827.59'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.52%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.58'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.34%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '42.53'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('E44').Value = '  -3.18%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '25.69'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +7.87%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.763'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '311.29'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.17%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('E51').Value = '  -1.41%  '
